# Generate Report for Handback
# Updates the "Overview", "zh-cn" and "de-de" sheets with a second
# handback entry (a0c7e227-...) and refreshes the timestamps / hashes
# of the existing (1162b4ef-... née af782309-...) entry.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Refresh row 2 (existing file) - file got renamed + re-timestamped
$ov.Range("A2").Value = "1162b4ef-c971-46e6-be00-412311172999.md"
$ov.Range("B2").Value = "e2e\1162b4ef-c971-46e6-be00-412311172999.md"
$ov.Range("G2").Value = "2016-09-02 15:17:11"
$ov.Range("G2").NumberFormat = $dateFmt

# New row 3 (second handback file)
$ov.Range("A3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-09-02 15:17:11"
$ov.Range("G3").NumberFormat = $dateFmt

$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b799d87378b4ed77b92f81499b8589a8383db71/e2e/a0c7e227-c2e6-4ead-915e-7a9f373574a4.md", [Type]::Missing, [Type]::Missing, "e2e\a0c7e227-c2e6-4ead-915e-7a9f373574a4.md") | Out-Null

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G3"))

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Refresh row 2 (existing file)
$zh.Range("A2").Value = "1162b4ef-c971-46e6-be00-412311172999.md"
$zh.Range("G2").Value = "1162b4ef-c971-46e6-be00-412311172999.e2395885892eadaa2ca0c9c596e451e6b5ef1111.zh-cn.xlf"
$zh.Range("H2").Value = "2016-09-02 15:16:58"
$zh.Range("H2").NumberFormat = $dateFmt
$zh.Range("I2").Value = "1162b4ef-c971-46e6-be00-412311172999.md"
$zh.Range("J2").Value = "1162b4ef-c971-46e6-be00-412311172999.e2395885892eadaa2ca0c9c596e451e6b5ef1111.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-02 15:17:32"
$zh.Range("K2").NumberFormat = $dateFmt

# New row 3 (second handback file)
$zh.Range("A3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "'True"
$zh.Range("G3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.5768b5312287569130b6c2c1c87346f96fc46b4b.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-02 15:16:58"
$zh.Range("H3").NumberFormat = $dateFmt
$zh.Range("I3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md"
$zh.Range("J3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.5768b5312287569130b6c2c1c87346f96fc46b4b.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-02 15:17:32"
$zh.Range("K3").NumberFormat = $dateFmt
$zh.Range("M3").Value = "'True"
$zh.Range("O3").Value = "'False"

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b799d87378b4ed77b92f81499b8589a8383db71/e2e/a0c7e227-c2e6-4ead-915e-7a9f373574a4.md", [Type]::Missing, [Type]::Missing, "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e11dd8afbf2f3b7360182265ea53b996a5b3cb20/e2e/a0c7e227-c2e6-4ead-915e-7a9f373574a4.md", [Type]::Missing, [Type]::Missing, "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md") | Out-Null

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P3"))

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Refresh row 2 (existing file)
$de.Range("A2").Value = "1162b4ef-c971-46e6-be00-412311172999.md"
$de.Range("G2").Value = "1162b4ef-c971-46e6-be00-412311172999.e2395885892eadaa2ca0c9c596e451e6b5ef1111.de-de.xlf"
$de.Range("H2").Value = "2016-09-02 15:17:11"
$de.Range("H2").NumberFormat = $dateFmt
$de.Range("I2").Value = "1162b4ef-c971-46e6-be00-412311172999.md"
$de.Range("J2").Value = "1162b4ef-c971-46e6-be00-412311172999.e2395885892eadaa2ca0c9c596e451e6b5ef1111.de-de.xlf"
$de.Range("K2").Value = "2016-09-02 15:17:39"
$de.Range("K2").NumberFormat = $dateFmt

# New row 3 (second handback file)
$de.Range("A3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "'True"
$de.Range("G3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.5768b5312287569130b6c2c1c87346f96fc46b4b.de-de.xlf"
$de.Range("H3").Value = "2016-09-02 15:17:11"
$de.Range("H3").NumberFormat = $dateFmt
$de.Range("I3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md"
$de.Range("J3").Value = "a0c7e227-c2e6-4ead-915e-7a9f373574a4.5768b5312287569130b6c2c1c87346f96fc46b4b.de-de.xlf"
$de.Range("K3").Value = "2016-09-02 15:17:39"
$de.Range("K3").NumberFormat = $dateFmt
$de.Range("M3").Value = "'True"
$de.Range("O3").Value = "'False"

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b799d87378b4ed77b92f81499b8589a8383db71/e2e/a0c7e227-c2e6-4ead-915e-7a9f373574a4.md", [Type]::Missing, [Type]::Missing, "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c4cd5cba62a8ab11d04a2548818a1aeac89086f7/e2e/a0c7e227-c2e6-4ead-915e-7a9f373574a4.md", [Type]::Missing, [Type]::Missing, "a0c7e227-c2e6-4ead-915e-7a9f373574a4.md") | Out-Null

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P3"))

Write-Host "Handback report rows added."
